$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.125.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.748.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5290"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2810"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06181"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.745.68"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07177"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6452"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.632"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.60"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.98%  "

$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.034.87"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006737"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.971.79"
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.321"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.722"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.226"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.52%  "

$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.810"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.91"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08305"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.803"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.644"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04630"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.649"
$ws.Range("D34").ClearFormats()

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.009"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6336"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01623"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.976"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9997"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.41"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3926"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7530"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.055"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.98%  "

$ws.Range("E45").Value = "  +3.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.356"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05348"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.59"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.04"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3482"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.572"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.33%  "
